# wireframes.docx: "Version 1." -> "Version 2."
#
# The target OOXML doesn't just swap the digit - it also re-shapes the run
# layout:
#   - "Version" (one run) becomes "Versi" + "on" (two runs)
#   - " 1." (one run) becomes " 2" (one run, trailing period removed)
#   - a brand new "." run is appended after the _GoBack bookmark
#
# Plain Range.Text edits never introduce a new run boundary: Word always
# coalesces adjacent runs that carry identical formatting. The trick used
# below is to drop a zero-width bookmark exactly at the desired split point
# and delete it immediately - that forces the engine to seat the text on
# either side in two separate <w:r> elements, and (unlike toggling a
# character property such as Bold on/off) it leaves no leftover <w:rPr/>
# behind once the bookmark is gone.

$d = $word.ActiveDocument

$para = $d.Paragraphs(1).Range
$paraStart = $para.Start

# --- 1. Split "Version" into "Versi" | "on" --------------------------------
$splitPoint = $d.Range($paraStart + 5, $paraStart + 5)
$d.Bookmarks.Add("__tmp_split", $splitPoint)
$d.Bookmarks("__tmp_split").Delete()

# --- 2. "1" -> "2" -----------------------------------------------------------
$digit = $d.Range($paraStart + 8, $paraStart + 9)
$digit.Text = "2"

# --- 3. Drop the old trailing "." (sat right after the "1"/"2") ------------
$oldPeriod = $d.Range($paraStart + 9, $paraStart + 10)
$oldPeriod.Delete()

# --- 4. Re-add "." as its own run, after the _GoBack bookmark --------------
$tail = $d.Range($paraStart + 9, $paraStart + 9)
$tail.InsertAfter(".")
